$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 59002.5
$ws.Range("J10").Value = 59002.5
$ws.Range("L10").Value = 59002.5
$ws.Range("N10").Value = -59588.5

$ws.Range("H13").Value = 8000
$ws.Range("J13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("N13").Value = -8338

$ws.Range("H53").Value = 18519606
$ws.Range("I53").Value = 47619856
$ws.Range("J53").Value = 1265.2727
$ws.Range("K53").Value = 47619856
$ws.Range("L53").Value = 1265.2727
$ws.Range("M53").Value = -47619219
$ws.Range("N53").Value = -2539.2727

$ws.Range("H80").Value = 12777.556
$ws.Range("I80").Value = 21999.6
$ws.Range("J80").Value = 1250
$ws.Range("K80").Value = 65998.79999999999
$ws.Range("L80").Value = 3750
$ws.Range("M80").Value = -65000.79999999999
$ws.Range("N80").Value = -5746

$ws.Range("H83").Value = 12777.556
$ws.Range("I83").Value = 21999.6
$ws.Range("J83").Value = 1250
$ws.Range("K83").Value = 197996.4
$ws.Range("L83").Value = 11250
$ws.Range("M83").Value = -193004.4
$ws.Range("N83").Value = -21234

$ws.Range("H88").Value = 1739.1818
$ws.Range("J88").Value = 1732.9
$ws.Range("L88").Value = 1732.9
$ws.Range("N88").Value = -2544.9

$ws.Range("H91").Value = 1739.1818
$ws.Range("J91").Value = 1732.9
$ws.Range("L91").Value = 1732.9
$ws.Range("N91").Value = -4540.9

$ws.Range("H100").Value = 9196.5
$ws.Range("I100").Value = 1249.5714
$ws.Range("K100").Value = 1249.5714
$ws.Range("M100").Value = -708.5714

$ws.Range("H125").Value = 12349860
$ws.Range("I125").Value = 2200
$ws.Range("J125").Value = 13893317
$ws.Range("K125").Value = 19800
$ws.Range("L125").Value = 125039853
$ws.Range("N125").Value = -125044773

$ws.Range("H129").Value = 1740.36
$ws.Range("I129").Value = 509.36365
$ws.Range("K129").Value = 1528.09095
$ws.Range("M129").Value = 3471.90905

$ws.Range("H132").Value = 2792.4167
$ws.Range("I132").Value = 2979
$ws.Range("K132").Value = 8937
$ws.Range("M132").Value = -6407

$ws.Range("H135").Value = 541963.4
$ws.Range("I135").Value = 690930.3
$ws.Range("K135").Value = 6218372.7
$ws.Range("M135").Value = -6215837.7

$ws.Range("H137").Value = 3422.4688
$ws.Range("I137").Value = 3917.1667
$ws.Range("J137").Value = 1938.375
$ws.Range("K137").Value = 11751.5001
$ws.Range("L137").Value = 5815.125
$ws.Range("M137").Value = -9201.500100000001
$ws.Range("N137").Value = -10915.125

$ws.Range("H138").Value = 4916.2144
$ws.Range("J138").Value = 7017.8
$ws.Range("L138").Value = 21053.4
$ws.Range("N138").Value = -31333.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 2646
$ws.Range("I31").Value = 2646
$ws.Range("K31").Value = 2646
$ws.Range("M31").Value = -2352

$ws.Range("H45").Value = 2826.4
$ws.Range("I45").Value = 2337.8
$ws.Range("K45").Value = 2337.8
$ws.Range("M45").Value = -1960.8

$ws.Range("H61").Value = 3800.6875
$ws.Range("I61").Value = 3089.3635
$ws.Range("J61").Value = 5365.6
$ws.Range("K61").Value = 3089.3635
$ws.Range("L61").Value = 5365.6
$ws.Range("M61").Value = -2877.3635
$ws.Range("N61").Value = -5789.6

$ws.Range("H110").Value = 126708.42
$ws.Range("I110").Value = 133308.88
$ws.Range("K110").Value = 133308.88
$ws.Range("M110").Value = -131263.88

$ws.Range("H132").Value = 3389.258
$ws.Range("I132").Value = 3385.5667
$ws.Range("K132").Value = 10156.7001
$ws.Range("M132").Value = -7626.7001

$ws.Range("H136").Value = 3800.6875
$ws.Range("I136").Value = 3089.3635
$ws.Range("J136").Value = 5365.6
$ws.Range("K136").Value = 9268.0905
$ws.Range("L136").Value = 16096.8
$ws.Range("M136").Value = -6718.0905
$ws.Range("N136").Value = -21196.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3378.4
$ws.Range("I99").Value = 3378.4
$ws.Range("K99").Value = 3378.4
$ws.Range("M99").Value = -1880.4

$ws.Range("H107").Value = 627338.5600000001
$ws.Range("I107").Value = 1622
$ws.Range("J107").Value = 2504488.2
$ws.Range("K107").Value = 1622
$ws.Range("L107").Value = 2504488.2
$ws.Range("M107").Value = 298
$ws.Range("N107").Value = -2508328.2

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

$ws.Range("H138").Value = 183333.33
$ws.Range("J138").Value = 183333.33
$ws.Range("L138").Value = 183333.33
$ws.Range("N138").Value = -193613.33

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 240321.53
$ws.Range("I134").Value = 2348.8206
$ws.Range("J134").Value = 3333966.8
$ws.Range("K134").Value = 7046.4618
$ws.Range("L134").Value = 10001900.4
$ws.Range("M134").Value = -4511.4618
$ws.Range("N134").Value = -10006970.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5805.4287
$ws.Range("I56").Value = 5805.4287
$ws.Range("K56").Value = 5805.4287
$ws.Range("M56").Value = -5275.4287

$ws.Range("H103").Value = 3498.3
$ws.Range("J103").Value = 5712.1665
$ws.Range("L103").Value = 17136.4995
$ws.Range("N103").Value = -18894.4995

$ws.Range("H134").Value = 1985.8889
$ws.Range("I134").Value = 1985.8889
$ws.Range("K134").Value = 5957.6667
$ws.Range("M134").Value = -887.6666999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1292.0625
$ws.Range("I107").Value = 796.5
$ws.Range("J107").Value = 1787.625
$ws.Range("K107").Value = 796.5
$ws.Range("L107").Value = 1787.625
$ws.Range("M107").Value = 1123.5
$ws.Range("N107").Value = -5627.625

$ws.Range("H132").Value = 98509
$ws.Range("I132").Value = 10657.286
$ws.Range("K132").Value = 31971.858
$ws.Range("M132").Value = -29441.858

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1429.6428
$ws.Range("I55").Value = 350.875
$ws.Range("K55").Value = 350.875
$ws.Range("M55").Value = -177.875

$ws.Range("H136").Value = 337897.5
$ws.Range("I136").Value = 439015.12
$ws.Range("K136").Value = 1317045.36
$ws.Range("M136").Value = -1314495.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 13692.385
$ws.Range("I74").Value = 4998
$ws.Range("J74").Value = 15273.182
$ws.Range("K74").Value = 4998
$ws.Range("L74").Value = 15273.182
$ws.Range("M74").Value = -4062
$ws.Range("N74").Value = -17145.182

$ws.Range("H77").Value = 13692.385
$ws.Range("I77").Value = 4998
$ws.Range("J77").Value = 15273.182
$ws.Range("K77").Value = 14994
$ws.Range("L77").Value = 45819.546
$ws.Range("M77").Value = -10314
$ws.Range("N77").Value = -55179.546

$ws.Range("H88").Value = 26275.857
$ws.Range("I88").Value = 13484.5
$ws.Range("J88").Value = 31392.4
$ws.Range("K88").Value = 13484.5
$ws.Range("L88").Value = 31392.4
$ws.Range("M88").Value = -13078.5
$ws.Range("N88").Value = -32204.4

$ws.Range("H91").Value = 26275.857
$ws.Range("I91").Value = 13484.5
$ws.Range("J91").Value = 31392.4
$ws.Range("K91").Value = 13484.5
$ws.Range("L91").Value = 31392.4
$ws.Range("M91").Value = -12080.5
$ws.Range("N91").Value = -34200.4

$ws.Range("H132").Value = 36663.234
$ws.Range("I132").Value = 2350.0833
$ws.Range("K132").Value = 7050.249899999999
$ws.Range("M132").Value = -4520.249899999999
